$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- New values for column F (climate change impact), keyed by row number ---
$newF = @{
  2 = [double]"-6.2580095e-06"
  3 = [double]"-7.5519471e-06"
  4 = [double]"-1.0296016e-05"
  5 = [double]"-9.4502253e-05"
  6 = [double]"-3.9706357e-07"
  7 = [double]"3.2189724e-06"
  8 = [double]"4.4189278e-06"
  9 = [double]"6.913417e-06"
  10 = [double]"6.6471273e-07"
  11 = [double]"2.9364015e-06"
  12 = [double]"8.898194300000001e-06"
  13 = [double]"2.9364015e-06"
  14 = [double]"-4.043606e-05"
  15 = [double]"-3.4915856e-05"
  16 = [double]"-1.7492004e-05"
  17 = [double]"-3.4075332e-05"
  18 = [double]"8.898194300000001e-06"
  19 = [double]"1.2974157e-05"
  20 = [double]"-4.1117567e-05"
  21 = [double]"-3.5438345e-05"
  22 = [double]"-1.7719172e-05"
  23 = [double]"7.8951352e-07"
  24 = [double]"-0.00010131732"
  25 = [double]"-3.3848163e-05"
  26 = [double]"2.703956e-06"
  27 = [double]"-6.5333769e-05"
  28 = [double]"3.5363848e-06"
  29 = [double]"1.3597169e-05"
  30 = [double]"-1.3765788e-05"
  31 = [double]"1.1102663e-05"
  32 = [double]"-5.9745415e-05"
  33 = [double]"-3.9754554e-05"
  34 = [double]"4.3224493e-06"
  35 = [double]"4.867649e-06"
  36 = [double]"-9.4268621e-06"
  37 = [double]"-1.16756e-06"
  38 = [double]"9.316420200000001e-06"
  39 = [double]"-5.0238422e-06"
  40 = [double]"8.898194300000001e-06"
  41 = [double]"-4.4001225e-07"
  42 = [double]"1.7476513e-06"
  43 = [double]"-3.6574189e-05"
  44 = [double]"-6.065409e-05"
  45 = [double]"-5.9972584e-05"
  46 = [double]"-5.7246557e-05"
  47 = [double]"1.9586504e-06"
  48 = [double]"-9.700111100000001e-05"
  49 = [double]"-1.6779238e-06"
  50 = [double]"-4.9977153e-05"
  51 = [double]"1.3294255e-06"
  52 = [double]"4.9747087e-06"
  53 = [double]"-1.3461915e-06"
  54 = [double]"6.0694765e-06"
  55 = [double]"-6.2188051e-06"
  56 = [double]"-6.7453388e-06"
  57 = [double]"-3.3669694e-07"
  58 = [double]"-1.2156939e-05"
  59 = [double]"1.0057781e-06"
  60 = [double]"1.0002919e-05"
  61 = [double]"3.4137737e-06"
  62 = [double]"3.7710089e-06"
  63 = [double]"-3.6574189e-05"
  64 = [double]"8.9047302e-06"
  65 = [double]"1.6519245e-05"
  66 = [double]"1.7746575e-05"
  67 = [double]"6.7198638e-06"
  68 = [double]"4.4263886e-06"
  69 = [double]"7.3783359e-06"
  70 = [double]"8.8247503e-06"
  71 = [double]"-2.8305242e-05"
  72 = [double]"-2.4534239e-05"
  73 = [double]"-1.2267119e-05"
  74 = [double]"-2.3852732e-05"
  75 = [double]"-2.8850447e-05"
  76 = [double]"-2.4761408e-05"
  77 = [double]"-1.2403421e-05"
  78 = [double]"3.9834556e-05"
  79 = [double]"2.4933202e-05"
  80 = [double]"3.4978742e-05"
  81 = [double]"-5.1976239e-05"
  82 = [double]"2.4750531e-05"
  83 = [double]"3.5522387e-05"
  84 = [double]"4.571823e-05"
  85 = [double]"2.776752e-05"
  86 = [double]"1.8724638e-05"
  87 = [double]"1.8724638e-05"
  88 = [double]"3.8657317e-05"
  89 = [double]"1.8724638e-05"
  90 = [double]"-2.2239833e-05"
  91 = [double]"-1.919577e-05"
  92 = [double]"-9.631960400000001e-06"
  93 = [double]"-1.8741432e-05"
  94 = [double]"3.8657317e-05"
  95 = [double]"2.9671069e-05"
  96 = [double]"-2.2716888e-05"
  97 = [double]"-1.9536523e-05"
  98 = [double]"-9.768261700000001e-06"
  99 = [double]"0"
  100 = [double]"3.9703542e-05"
  101 = [double]"-5.5724526e-05"
  102 = [double]"-1.861649e-05"
  103 = [double]"4.0891044e-05"
  104 = [double]"-3.5938116e-05"
  105 = [double]"3.2932038e-05"
  106 = [double]"4.4378552e-05"
  107 = [double]"3.1827006e-05"
  108 = [double]"3.4841811e-05"
  109 = [double]"-3.2939487e-05"
  110 = [double]"-2.1865004e-05"
  111 = [double]"2.9969816e-05"
  112 = [double]"2.1110224e-05"
  113 = [double]"3.4212279e-05"
  114 = [double]"4.1767358e-05"
  115 = [double]"1.617692e-05"
  116 = [double]"1.8692589e-05"
  117 = [double]"3.8657317e-05"
  118 = [double]"1.7449537e-05"
  119 = [double]"1.2379155e-05"
  120 = [double]"-2.0104446e-05"
  121 = [double]"-3.3166656e-05"
  122 = [double]"-3.2939487e-05"
  123 = [double]"-3.1349305e-05"
  124 = [double]"3.239928e-05"
  125 = [double]"-5.3339252e-05"
  126 = [double]"2.201579e-05"
  127 = [double]"-2.7487434e-05"
  128 = [double]"3.4950419e-05"
  129 = [double]"3.8005063e-05"
  130 = [double]"2.8004027e-05"
  131 = [double]"3.9599603e-05"
  132 = [double]"3.0582553e-05"
  133 = [double]"3.0669375e-05"
  134 = [double]"1.7745946e-05"
  135 = [double]"2.3962912e-05"
  136 = [double]"1.8997553e-05"
  137 = [double]"3.8626198e-05"
  138 = [double]"2.9470044e-05"
  139 = [double]"1.1222148e-05"
  140 = [double]"-2.0104446e-05"
  141 = [double]"4.0163168e-05"
  142 = [double]"4.5346975e-05"
  143 = [double]"3.9622938e-05"
  144 = [double]"3.9318598e-05"
  145 = [double]"3.3049667e-05"
  146 = [double]"2.9186548e-05"
  147 = [double]"3.3449857e-05"
  148 = [double]"-5.1476468e-05"
  149 = [double]"2.9992532e-05"
  150 = [double]"-4.7614597e-05"
  151 = [double]"9.6195154e-05"
  152 = [double]"6.4687756e-05"
  153 = [double]"9.0339797e-05"
  154 = [double]"0"
  155 = [double]"5.548648e-05"
  156 = [double]"7.500433799999999e-05"
  157 = [double]"9.6195154e-05"
  158 = [double]"5.3255868e-05"
  159 = [double]"8.8945664e-05"
  160 = [double]"3.7920409e-05"
  161 = [double]"0.00015000868"
  162 = [double]"3.7920409e-05"
  163 = [double]"4.7958164e-05"
  164 = [double]"8.7272705e-05"
  165 = [double]"0"
  166 = [double]"0"
  167 = [double]"8.7551532e-05"
  168 = [double]"0"
  169 = [double]"6.8870154e-05"
  170 = [double]"8.1975001e-05"
  171 = [double]"8.7551532e-05"
  172 = [double]"6.3851276e-05"
  173 = [double]"0"
  174 = [double]"0"
  175 = [double]"6.134183799999999e-05"
  176 = [double]"4.0987501e-05"
  177 = [double]"8.7551532e-05"
  178 = [double]"9.4243369e-05"
  179 = [double]"2.4536735e-05"
  180 = [double]"4.7679337e-05"
  181 = [double]"7.500433799999999e-05"
  182 = [double]"2.5373215e-05"
  183 = [double]"0"
  184 = [double]"0"
  185 = [double]"0"
  186 = [double]"0"
  187 = [double]"6.9427807e-05"
  188 = [double]"0"
  189 = [double]"6.3851276e-05"
  190 = [double]"0"
  191 = [double]"7.584081699999999e-05"
  192 = [double]"7.8350256e-05"
  193 = [double]"6.3851276e-05"
  194 = [double]"8.058086800000001e-05"
  195 = [double]"7.5561991e-05"
  196 = [double]"7.639847e-05"
  197 = [double]"3.9872194e-05"
  198 = [double]"6.803367399999999e-05"
  199 = [double]"4.0987501e-05"
  200 = [double]"7.3610205e-05"
  201 = [double]"6.134183799999999e-05"
  202 = [double]"2.0354337e-05"
  203 = [double]"0"
  204 = [double]"7.8350256e-05"
  205 = [double]"8.058086800000001e-05"
  206 = [double]"6.6360715e-05"
  207 = [double]"7.9186736e-05"
  208 = [double]"6.803367399999999e-05"
  209 = [double]"5.5765307e-05"
  210 = [double]"6.357245e-05"
  211 = [double]"0"
  212 = [double]"6.134183799999999e-05"
  213 = [double]"0"
  214 = [double]"1.5057262e-06"
  215 = [double]"1.5057262e-06"
  216 = [double]"2.3570371e-06"
  217 = [double]"1.9804727e-06"
  218 = [double]"2.3784905e-06"
  219 = [double]"1.9706845e-06"
  220 = [double]"2.6394668e-06"
  221 = [double]"-1.5629219e-06"
  222 = [double]"1.5521745e-06"
  223 = [double]"1.6452353e-06"
  224 = [double]"2.188438e-06"
  225 = [double]"1.1672654e-06"
  226 = [double]"0"
  227 = [double]"-8.9488573e-05"
  228 = [double]"-3.7428693e-05"
  229 = [double]"-1.0492728"
  230 = [double]"-8.9543528e-06"
  231 = [double]"-0.00065208838"
  232 = [double]"-0.00023101208"
  233 = [double]"-0.59344072"
  234 = [double]"-0.79457065"
  235 = [double]"-0.80671512"
  236 = [double]"-0.0054581917"
  237 = [double]"-2.5666468e-05"
  238 = [double]"-0.00022081554"
  239 = [double]"-5.5512477e-05"
  240 = [double]"0"
  241 = [double]"-5.9294521e-06"
  242 = [double]"1.5664638e-06"
  243 = [double]"0"
  244 = [double]"8.4872315e-06"
  245 = [double]"-7.5568655e-06"
  246 = [double]"0"
  247 = [double]"2.5248049e-06"
  248 = [double]"-9.8143876e-07"
  249 = [double]"-7.5568655e-06"
  250 = [double]"-2.5898212e-05"
  251 = [double]"0"
  252 = [double]"0"
  253 = [double]"1.0608683e-05"
  254 = [double]"2.3274797e-05"
  255 = [double]"-7.5568656e-06"
  256 = [double]"-6.5610365e-06"
  257 = [double]"2.5555444e-05"
  258 = [double]"1.423287e-05"
  259 = [double]"1.1684532e-06"
  260 = [double]"3.2939478e-06"
  261 = [double]"0"
  262 = [double]"0"
  263 = [double]"0"
  264 = [double]"-7.5568655e-06"
  265 = [double]"1.5664638e-06"
  266 = [double]"-2.5898212e-05"
  267 = [double]"0"
  268 = [double]"-2.0170337e-06"
  269 = [double]"-6.6081168e-06"
  270 = [double]"-5.6994412e-06"
  271 = [double]"0"
  272 = [double]"-4.7091004e-07"
  273 = [double]"0"
  274 = [double]"6.2485026e-05"
  275 = [double]"3.6347528e-06"
  276 = [double]"-4.3230187e-05"
  277 = [double]"0"
  278 = [double]"2.483395e-06"
  279 = [double]"3.6761444e-07"
  280 = [double]"1.2416975e-06"
  281 = [double]"1.2416975e-06"
  282 = [double]"8.2779832e-07"
  283 = [double]"2.483395e-06"
  284 = [double]"2.0694958e-06"
  285 = [double]"2.483395e-06"
  286 = [double]"4.2102807e-05"
  287 = [double]"-8.7551532e-05"
  288 = [double]"-8.7551532e-05"
  289 = [double]"-1.2212602e-05"
  290 = [double]"-5.3813521e-05"
  291 = [double]"-6.8591327e-05"
  292 = [double]"0"
  293 = [double]"-1.3941327e-06"
  294 = [double]"-3.3459184e-06"
}

# --- Shift data columns: new D = old E, new E = old F, new F = newly computed value ---
for ($i = 2; $i -le 294; $i++) {
    $oldE = $ws.Cells.Item($i, 5).Value2
    $oldF = $ws.Cells.Item($i, 6).Value2
    $ws.Cells.Item($i, 4).Value = $oldE
    $ws.Cells.Item($i, 5).Value = $oldF
    $ws.Cells.Item($i, 6).Value = $newF[$i]
}

# --- Add header comments describing each column's data type ---
$cA = $ws.Range("A1").AddComment("Data type: Categorical (text)")
$cB = $ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)")
$cC = $ws.Range("C1").AddComment("Data type: Categorical (text)")
$cD = $ws.Range("D1").AddComment("Data type: Carbon footprint")
$cE = $ws.Range("E1").AddComment("Data type: Cumulative energy demand")
$cF = $ws.Range("F1").AddComment("Data type: Climate change impact")
$cG = $ws.Range("G1").AddComment("Data type: Categorical (text)")
